# nomina: actualizacion de plantilla para subir banco y cuenta banco
#
# Insert two new columns ("BANCO" and "CUENTA_BANCO") right after the
# existing "NUMERO_HIJOS" column (col W) and before "NIVEL_FORMACION"
# (old col X), shifting everything from old X.. onward two columns to
# the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank columns at X:Y (pushes old X:AE -> Z:AG), carrying over
# formatting (header style, etc.) from the column that used to sit there.
$ws.Range("X1:Y1").EntireColumn.Insert()

# New header labels for the freshly inserted columns.
$ws.Range("X1").Value = "BANCO"
$ws.Range("Y1").Value = "CUENTA_BANCO"

# Match the column widths Excel computed for the new headers (best-fit).
$ws.Columns.Item(24).ColumnWidth = 6.25
$ws.Columns.Item(25).ColumnWidth = 15

# Restore the view: scrolled so column M is left-most visible, and the
# active cell/selection on X6.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 13
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("X6").Select()
